$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174654245376587
$ws.Range("B1").Value = 2.286669969558716
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.141146898269653
$ws.Range("E1").Value = 1.045525312423706
